# Apply the 3 Feb 2023 19:08 UTC GitHub Actions symbol-list refresh:
# updated Price (D) / Volume(1h) (E) figures and bumped the Hora (G)
# column from "18" to "19" for every data row (2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (matches the sheet's existing
# text-typed Price/Volume/Hora columns) instead of letting Excel
# auto-convert numeric-looking / percent-looking strings to numbers.
function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

# Column G ("Hora"): every row 2-51 moves from "18" to "19".
Set-TextValue $ws.Range("G2:G51") "19"

# Row 2
Set-TextValue $ws.Range("D2") "332.82"
Set-TextValue $ws.Range("E2") "0.91%"
# Row 3
Set-TextValue $ws.Range("D3") "41.25"
Set-TextValue $ws.Range("E3") "1.20%"
# Row 4
Set-TextValue $ws.Range("D4") "5.683"
Set-TextValue $ws.Range("E4") "-7.03%"
# Row 5
Set-TextValue $ws.Range("D5") "0.08080"
Set-TextValue $ws.Range("E5") "-1.32%"
# Row 6
Set-TextValue $ws.Range("D6") "2.039"
Set-TextValue $ws.Range("E6") "3.48%"
# Row 7
Set-TextValue $ws.Range("D7") "8.748"
Set-TextValue $ws.Range("E7") "-0.56%"
# Row 8
Set-TextValue $ws.Range("D8") "4.544"
Set-TextValue $ws.Range("E8") "-1.01%"
# Row 10
Set-TextValue $ws.Range("D10") "0.9227"
Set-TextValue $ws.Range("E10") "-2.98%"
# Row 11
Set-TextValue $ws.Range("D11") "0.1257"
Set-TextValue $ws.Range("E11") "-7.50%"
# Row 12
Set-TextValue $ws.Range("D12") "0.1944"
Set-TextValue $ws.Range("E12") "-3.18%"
# Row 13
Set-TextValue $ws.Range("D13") "8.773"
Set-TextValue $ws.Range("E13") "-16.00%"
# Row 14
Set-TextValue $ws.Range("D14") "0.09512"
Set-TextValue $ws.Range("E14") "2.66%"
# Row 15
Set-TextValue $ws.Range("D15") "0.03706"
Set-TextValue $ws.Range("E15") "4.26%"
# Row 16
Set-TextValue $ws.Range("D16") "0.1051"
Set-TextValue $ws.Range("E16") "9.07%"
# Row 17
Set-TextValue $ws.Range("D17") "0.001309"
Set-TextValue $ws.Range("E17") "-0.28%"
# Row 18
Set-TextValue $ws.Range("D18") "0.006247"
Set-TextValue $ws.Range("E18") "-2.75%"
# Row 19
Set-TextValue $ws.Range("E19") "0.27%"
# Row 20
Set-TextValue $ws.Range("E20") "-1.26%"
# Row 21
Set-TextValue $ws.Range("D21") "0.1420"
Set-TextValue $ws.Range("E21") "-1.25%"
# Row 22
Set-TextValue $ws.Range("D22") "0.2656"
Set-TextValue $ws.Range("E22") "9.20%"
# Row 23
Set-TextValue $ws.Range("D23") "0.04434"
# Row 24
Set-TextValue $ws.Range("D24") "0.001261"
Set-TextValue $ws.Range("E24") "0.09%"
# Row 25
Set-TextValue $ws.Range("D25") "0.004300"
Set-TextValue $ws.Range("E25") "-2.38%"
# Row 26
Set-TextValue $ws.Range("E26") "13.98%"
# Row 39
Set-TextValue $ws.Range("D39") "0.02874"
Set-TextValue $ws.Range("E39") "13.97%"
# Row 40
Set-TextValue $ws.Range("D40") "0.05499"
Set-TextValue $ws.Range("E40") "3.89%"
# Row 41
Set-TextValue $ws.Range("D41") "0.007791"
Set-TextValue $ws.Range("E41") "4.19%"
# Row 42
Set-TextValue $ws.Range("D42") "0.009969"
Set-TextValue $ws.Range("E42") "10.55%"
# Row 43
Set-TextValue $ws.Range("E43") "-2.12%"
# Row 44
Set-TextValue $ws.Range("D44") "0.002134"
Set-TextValue $ws.Range("E44") "4.02%"
# Row 45
Set-TextValue $ws.Range("D45") "0.01110"
Set-TextValue $ws.Range("E45") "4.98%"
# Row 46
Set-TextValue $ws.Range("D46") "0.00006823"
Set-TextValue $ws.Range("E46") "0.59%"
# Row 47
Set-TextValue $ws.Range("D47") "0.00000000752"
Set-TextValue $ws.Range("E47") "0.23%"
# Row 48
Set-TextValue $ws.Range("D48") "0.002284"
Set-TextValue $ws.Range("E48") "26.90%"
# Row 49
Set-TextValue $ws.Range("D49") "0.003024"
Set-TextValue $ws.Range("E49") "-13.48%"
# Row 50
Set-TextValue $ws.Range("D50") "0.00002105"
Set-TextValue $ws.Range("E50") "0.23%"
# Row 51
Set-TextValue $ws.Range("D51") "0.0002005"
Set-TextValue $ws.Range("E51") "0.23%"
